# Scheduled market-data refresh: updates the cached currentAveragePrice /
# LevePrice / LeveProfit columns (H:N) on each leve-items sheet with
# freshly pulled Universalis price data. Values are static snapshots
# (no formulas in these columns), so each changed cell is written directly.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 202.9
$ws.Range("I5").Value = 144.83333
$ws.Range("J5").Value = 290
$ws.Range("K5").Value = 144.83333
$ws.Range("L5").Value = 290
$ws.Range("M5").Value = -29.83332999999999
$ws.Range("N5").Value = -520
$ws.Range("H9").Value = 398.125
$ws.Range("I9").Value = 397.66666
$ws.Range("J9").Value = 399.5
$ws.Range("K9").Value = 397.66666
$ws.Range("L9").Value = 399.5
$ws.Range("M9").Value = -228.66666
$ws.Range("N9").Value = -737.5
$ws.Range("H40").Value = 4120.5713
$ws.Range("I40").Value = 4017.2856
$ws.Range("J40").Value = 4223.857
$ws.Range("K40").Value = 4017.2856
$ws.Range("L40").Value = 4223.857
$ws.Range("M40").Value = -3842.2856
$ws.Range("N40").Value = -4573.857
$ws.Range("H51").Value = 10416.258
$ws.Range("J51").Value = 7001
$ws.Range("L51").Value = 7001
$ws.Range("N51").Value = -7969
$ws.Range("H69").Value = 3660
$ws.Range("I69").Value = 2995
$ws.Range("J69").Value = 4990
$ws.Range("K69").Value = 8985
$ws.Range("L69").Value = 14970
$ws.Range("M69").Value = -8111
$ws.Range("N69").Value = -16718
$ws.Range("H70").Value = 577.2308
$ws.Range("I70").Value = 611.1111
$ws.Range("K70").Value = 1833.3333
$ws.Range("M70").Value = -1563.3333
$ws.Range("H72").Value = 3660
$ws.Range("I72").Value = 2995
$ws.Range("J72").Value = 4990
$ws.Range("K72").Value = 26955
$ws.Range("L72").Value = 44910
$ws.Range("M72").Value = -22587
$ws.Range("N72").Value = -53646
$ws.Range("H73").Value = 577.2308
$ws.Range("I73").Value = 611.1111
$ws.Range("K73").Value = 1833.3333
$ws.Range("M73").Value = -897.3332999999998
$ws.Range("H86").Value = 1709.1111
$ws.Range("I86").Value = 1297.7142
$ws.Range("J86").Value = 3149
$ws.Range("K86").Value = 1297.7142
$ws.Range("L86").Value = 3149
$ws.Range("M86").Value = -174.7141999999999
$ws.Range("N86").Value = -5395
$ws.Range("H87").Value = 30349
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 30349
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 30349
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -32845
$ws.Range("H89").Value = 1709.1111
$ws.Range("I89").Value = 1297.7142
$ws.Range("J89").Value = 3149
$ws.Range("K89").Value = 6488.571
$ws.Range("L89").Value = 15745
$ws.Range("M89").Value = -872.5709999999999
$ws.Range("N89").Value = -26977
$ws.Range("H90").Value = 30349
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 30349
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 91047
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -103527
$ws.Range("H129").Value = 1793.5
$ws.Range("I129").Value = 1306.625
$ws.Range("J129").Value = 1908.0588
$ws.Range("K129").Value = 3919.875
$ws.Range("L129").Value = 5724.1764
$ws.Range("M129").Value = 1080.125
$ws.Range("N129").Value = -15724.1764
$ws.Range("H131").Value = 3383.4167
$ws.Range("I131").Value = 2373.5
$ws.Range("J131").Value = 3888.375
$ws.Range("K131").Value = 7120.5
$ws.Range("L131").Value = 11665.125
$ws.Range("M131").Value = -2080.5
$ws.Range("N131").Value = -21745.125
$ws.Range("H132").Value = 127536.875
$ws.Range("I132").Value = 168540.83
$ws.Range("K132").Value = 505622.49
$ws.Range("M132").Value = -503092.49
$ws.Range("H137").Value = 13835.617
$ws.Range("I137").Value = 2750.3462
$ws.Range("K137").Value = 8251.0386
$ws.Range("M137").Value = -5701.0386
$ws.Range("H138").Value = 2145.3845
$ws.Range("I138").Value = 1801.6666
$ws.Range("J138").Value = 2918.75
$ws.Range("K138").Value = 5404.9998
$ws.Range("L138").Value = 8756.25
$ws.Range("M138").Value = -264.9997999999996
$ws.Range("N138").Value = -19036.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 737.04443
$ws.Range("I32").Value = 737.04443
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 737.04443
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -450.04443
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 3254.5454
$ws.Range("I61").Value = 2860.75
$ws.Range("J61").Value = 4304.6665
$ws.Range("K61").Value = 2860.75
$ws.Range("L61").Value = 4304.6665
$ws.Range("M61").Value = -2648.75
$ws.Range("N61").Value = -4728.6665
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H110").Value = 5830.5713
$ws.Range("I110").Value = 4514.222
$ws.Range("J110").Value = 8200
$ws.Range("K110").Value = 4514.222
$ws.Range("L110").Value = 8200
$ws.Range("M110").Value = -2469.222
$ws.Range("N110").Value = -12290
$ws.Range("H122").Value = 2348.1333
$ws.Range("I122").Value = 2444.3572
$ws.Range("K122").Value = 7333.071599999999
$ws.Range("M122").Value = -4883.071599999999
$ws.Range("H136").Value = 3254.5454
$ws.Range("I136").Value = 2860.75
$ws.Range("J136").Value = 4304.6665
$ws.Range("K136").Value = 8582.25
$ws.Range("L136").Value = 12913.9995
$ws.Range("M136").Value = -6032.25
$ws.Range("N136").Value = -18013.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 51000
$ws.Range("I9").Value = 51000
$ws.Range("K9").Value = 51000
$ws.Range("M9").Value = -50832
$ws.Range("H44").Value = 12500
$ws.Range("I44").Value = 12500
$ws.Range("K44").Value = 12500
$ws.Range("M44").Value = -12003
$ws.Range("H80").Value = 142.66667
$ws.Range("J80").Value = 123.333336
$ws.Range("L80").Value = 123.333336
$ws.Range("N80").Value = -2119.333336
$ws.Range("H83").Value = 142.66667
$ws.Range("J83").Value = 123.333336
$ws.Range("L83").Value = 616.66668
$ws.Range("N83").Value = -10600.66668
$ws.Range("H134").Value = 83347050
$ws.Range("I134").Value = 41682664
$ws.Range("K134").Value = 125047992
$ws.Range("M134").Value = -125045457

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1095.8235
$ws.Range("J22").Value = 4250
$ws.Range("L22").Value = 4250
$ws.Range("N22").Value = -4950
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
$ws.Range("H122").Value = 2356.2856
$ws.Range("I122").Value = 1997.5
$ws.Range("K122").Value = 5992.5
$ws.Range("M122").Value = -3542.5
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 3913.3333
$ws.Range("I61").Value = 3913.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 11739.9999
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -11524.9999
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 7358.6
$ws.Range("I53").Value = 7348
$ws.Range("K53").Value = 7348
$ws.Range("M53").Value = -6717
$ws.Range("H80").Value = 6015.5
$ws.Range("I80").Value = 2627
$ws.Range("K80").Value = 2627
$ws.Range("M80").Value = -1629
$ws.Range("H83").Value = 6015.5
$ws.Range("I83").Value = 2627
$ws.Range("K83").Value = 13135
$ws.Range("M83").Value = -8143
$ws.Range("H113").Value = 2494.2
$ws.Range("I113").Value = 1735.5
$ws.Range("K113").Value = 1735.5
$ws.Range("M113").Value = 434.5
$ws.Range("H126").Value = 4821.857
$ws.Range("I126").Value = 7067.6665
$ws.Range("J126").Value = 3137.5
$ws.Range("K126").Value = 21202.9995
$ws.Range("L126").Value = 9412.5
$ws.Range("M126").Value = -18732.9995
$ws.Range("N126").Value = -14352.5
$ws.Range("H132").Value = 2021.9412
$ws.Range("I132").Value = 1876.0769
$ws.Range("K132").Value = 5628.2307
$ws.Range("M132").Value = -3098.2307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3696.8462
$ws.Range("I40").Value = 3456.6
$ws.Range("J40").Value = 4497.6665
$ws.Range("K40").Value = 3456.6
$ws.Range("L40").Value = 4497.6665
$ws.Range("M40").Value = -3320.6
$ws.Range("N40").Value = -4769.6665
$ws.Range("H132").Value = 3939.4167
$ws.Range("I132").Value = 3920.9443
$ws.Range("K132").Value = 11762.8329
$ws.Range("M132").Value = -9232.832900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2385.3333
$ws.Range("I113").Value = 874.2
$ws.Range("K113").Value = 2622.6
$ws.Range("M113").Value = -452.6000000000004
$ws.Range("H122").Value = 4692
$ws.Range("I122").Value = 4659.5
$ws.Range("J122").Value = 4724.5
$ws.Range("K122").Value = 13978.5
$ws.Range("L122").Value = 14173.5
$ws.Range("M122").Value = -11528.5
$ws.Range("N122").Value = -19073.5
$ws.Range("H126").Value = 3293.25
$ws.Range("I126").Value = 2736.75
$ws.Range("K126").Value = 8210.25
$ws.Range("M126").Value = -5740.25
$ws.Range("H132").Value = 1707.439
$ws.Range("I132").Value = 1732.55
$ws.Range("K132").Value = 5197.65
$ws.Range("M132").Value = -2667.65
$ws.Range("H136").Value = 1261.5
$ws.Range("I136").Value = 1170.0526
$ws.Range("K136").Value = 3510.1578
$ws.Range("M136").Value = -960.1578
